$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Percentages (numeric, no shared-string impact)
$ws.Range("E7").Value = 0.5
$ws.Range("E8").Value = 0
$ws.Range("E9").Value = 0.5

# Text values - set in the order that reproduces the original shared-string table order
$ws.Range("F9").Value = "Responsible for implementing the data preprocessing mechanisms for the datasets utilized as well as to proper evaluate them and document their behaviour. In addition, he also recorded, edited and reviewed the video presentation."
$ws.Range("F8").Value = "Did not develop any kind of work"
$ws.Range("G8").Value = "The Student did not engage on the project's development throughout the pratical classes nor on extra efforts apart from university. Moreover, he did not display any effort or commitment on working in the project."
$ws.Range("G9").Value = "~"
$ws.Range("G7").Value = "~"
$ws.Range("F7").Value = "Developed the ID3 Algorithm as well as a few data visualization functions to better visualize the behaviour of the model on each dataset. In addition, he also integrated the ID3 Algorithm inside the Connect Four game and analysed the result's obtained."

# Row heights
$ws.Rows.Item(7).RowHeight = 60
$ws.Rows.Item(8).RowHeight = 80.1
$ws.Rows.Item(9).RowHeight = 60

# Column G width (target stored width 50.7109375; the host quantizes to pixel
# steps, so 49.8 is the input that lands closest on the achievable grid)
$ws.Columns.Item(7).ColumnWidth = 49.8

# G8 wraps text (matches the wrap style already used elsewhere)
$ws.Range("G8").WrapText = $true

# F7: text wraps & shrinks to fit, stored as text number format
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").WrapText = $true
$ws.Range("F7").ShrinkToFit = $true

# Partial borders - order matters for matching style indices:
# E9 loses its right edge, G9 loses its left edge, F8 loses its bottom edge
$ws.Range("E9").Borders.Item(10).LineStyle = 0
$ws.Range("G9").Borders.Item(7).LineStyle = 0
$ws.Range("F8").Borders.Item(9).LineStyle = 0

# New row 13 with an empty underlined cell at F13
$ws.Range("F13").Font.Underline = 1

# Selection
$ws.Range("F7").Select()
